$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Test", "da0cb5dd-123d-4b56-b46d-f6973aca7c00", "2025-03-d 21:12:19"),
    @("Test", "dd4e5bd3-76a0-4e9e-866b-5915b554b789", "2025-03-d 21:12:26"),
    @("Test", "c23811a0-d62c-4e4b-8cf5-840c9f527812", "2025-03-d 21:12:26"),
    @("Nguyen", "d89c668a-a383-4f5d-9a6f-dc9a4313989b", "2025-03-d 15:54:57"),
    @("PHONG_OLD", "63c2e3ba-ae2c-4ff7-a47f-9476fefddde5", "2025-03-d 15:56:56"),
    @("Test", "a03726e7-a734-41ae-9f7e-da5b7246a331", "2025-03-d 15:58:16"),
    @("Test", "598fec19-ac64-4230-92f0-63cd35d531e0", "2025-03-d 16:03:41"),
    @("Test", "3f358249-458d-4d40-a7ad-36127ca0a8e8", "2025-03-d 16:13:38"),
    @("Test", "80a11c55-9172-4fd3-9ba1-ef5755a2db80", "2025-03-d 16:25:06"),
    @("Test", "5081cb67-6626-4fb4-b464-1ea38cc43d65", "2025-03-d 16:30:44"),
    @("PHONG_OLD", "e8d307cd-ff42-4a49-a548-1b581edc8318", "2025-03-d 16:34:38"),
    @("TYLER", "5976e3ea-2dd9-49d2-9ccf-efcff2d56f06", "2025-03-d 16:36:35"),
    @("Test", "f107d882-0baf-44e0-b09f-5ad7a7fce4c7", "2025-03-d 16:39:45")
)

$startRow = 14
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
